$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 10).Value = 1.95
$ws.Cells.Item(2, 12).Value = 6.5
$ws.Cells.Item(2, 13).Value = 1.03
$ws.Cells.Item(2, 15).Value = 1.19
$ws.Cells.Item(2, 17).Value = 1.67
$ws.Cells.Item(2, 18).Value = 2.2
$ws.Cells.Item(2, 19).Value = 2.7
$ws.Cells.Item(2, 20).Value = 1.41
$ws.Cells.Item(2, 21).Value = 1.3
$ws.Cells.Item(2, 22).Value = 3.4
$ws.Cells.Item(2, 27).Value = 8.5
$ws.Cells.Item(2, 31).Value = 15
$ws.Cells.Item(2, 32).Value = 9

# Row 3
$ws.Cells.Item(3, 13).Value = 1.08
$ws.Cells.Item(3, 15).Value = 1.46
$ws.Cells.Item(3, 16).Value = 2.57
$ws.Cells.Item(3, 20).Value = 1.13

# Row 7
$ws.Cells.Item(7, 8).Value = 3.05
$ws.Cells.Item(7, 11).Value = 1.87
$ws.Cells.Item(7, 12).Value = 5
$ws.Cells.Item(7, 26).Value = 7.2
$ws.Cells.Item(7, 36).Value = 9.25
$ws.Cells.Item(7, 40).Value = 55

# Row 8
$ws.Cells.Item(8, 7).Value = 1.73
$ws.Cells.Item(8, 8).Value = 3.4
$ws.Cells.Item(8, 9).Value = 5.25
$ws.Cells.Item(8, 10).Value = 2.5
$ws.Cells.Item(8, 25).Value = 4.75
$ws.Cells.Item(8, 28).Value = 13
$ws.Cells.Item(8, 29).Value = 19
$ws.Cells.Item(8, 31).Value = 6
$ws.Cells.Item(8, 33).Value = 26
$ws.Cells.Item(8, 42).Value = 2
$ws.Cells.Item(8, 43).Value = 1.85

# Row 9
$ws.Cells.Item(9, 7).Value = 2.87
$ws.Cells.Item(9, 9).Value = 2.3
$ws.Cells.Item(9, 10).Value = 3.75
$ws.Cells.Item(9, 12).Value = 3.1
$ws.Cells.Item(9, 13).Value = 1.04
$ws.Cells.Item(9, 14).Value = 10
$ws.Cells.Item(9, 15).Value = 1.33
$ws.Cells.Item(9, 20).Value = 1.19
$ws.Cells.Item(9, 36).Value = 7
$ws.Cells.Item(9, 38).Value = 9.5
$ws.Cells.Item(9, 39).Value = 21

# Row 10
$ws.Cells.Item(10, 7).Value = 1.87
$ws.Cells.Item(10, 12).Value = 5
$ws.Cells.Item(10, 13).Value = 1.07
$ws.Cells.Item(10, 15).Value = 1.47
$ws.Cells.Item(10, 20).Value = 1.11
$ws.Cells.Item(10, 21).Value = 1.57
$ws.Cells.Item(10, 22).Value = 2.25
$ws.Cells.Item(10, 23).Value = 2.2
$ws.Cells.Item(10, 24).Value = 1.62
$ws.Cells.Item(10, 25).Value = 5.5
$ws.Cells.Item(10, 28).Value = 17
$ws.Cells.Item(10, 29).Value = 19
$ws.Cells.Item(10, 31).Value = 6.5
$ws.Cells.Item(10, 38).Value = 15
$ws.Cells.Item(10, 44).Value = 4.3
$ws.Cells.Item(10, 45).Value = 1.23

# Row 11
$ws.Cells.Item(11, 7).Value = 3.2
$ws.Cells.Item(11, 8).Value = 2.87
$ws.Cells.Item(11, 9).Value = 2.5
$ws.Cells.Item(11, 10).Value = 4
$ws.Cells.Item(11, 12).Value = 3.5
$ws.Cells.Item(11, 13).Value = 1.1
$ws.Cells.Item(11, 14).Value = 6
$ws.Cells.Item(11, 15).Value = 1.58
$ws.Cells.Item(11, 20).Value = 1.08
$ws.Cells.Item(11, 26).Value = 13
$ws.Cells.Item(11, 28).Value = 34
$ws.Cells.Item(11, 39).Value = 26
$ws.Cells.Item(11, 40).Value = 29
$ws.Cells.Item(11, 44).Value = 4.8
$ws.Cells.Item(11, 45).Value = 1.19

# Row 12
$ws.Cells.Item(12, 7).Value = 2.6
$ws.Cells.Item(12, 8).Value = 3.5
$ws.Cells.Item(12, 9).Value = 2.35
$ws.Cells.Item(12, 10).Value = 3.2
$ws.Cells.Item(12, 11).Value = 2.3
$ws.Cells.Item(12, 12).Value = 3
$ws.Cells.Item(12, 14).Value = 15
$ws.Cells.Item(12, 15).Value = 1.2
$ws.Cells.Item(12, 16).Value = 4.33
$ws.Cells.Item(12, 17).Value = 1.67
$ws.Cells.Item(12, 18).Value = 2.15
$ws.Cells.Item(12, 19).Value = 2.63
$ws.Cells.Item(12, 20).Value = 1.44
$ws.Cells.Item(12, 21).Value = 1.33
$ws.Cells.Item(12, 22).Value = 3.25
$ws.Cells.Item(12, 23).Value = 1.57
$ws.Cells.Item(12, 24).Value = 2.25
$ws.Cells.Item(12, 26).Value = 15
$ws.Cells.Item(12, 27).Value = 10
$ws.Cells.Item(12, 28).Value = 26
$ws.Cells.Item(12, 29).Value = 19
$ws.Cells.Item(12, 30).Value = 23
$ws.Cells.Item(12, 31).Value = 15
$ws.Cells.Item(12, 32).Value = 7
$ws.Cells.Item(12, 34).Value = 41
$ws.Cells.Item(12, 35).Value = 126
$ws.Cells.Item(12, 36).Value = 11
$ws.Cells.Item(12, 37).Value = 13
$ws.Cells.Item(12, 38).Value = 9.5
$ws.Cells.Item(12, 39).Value = 23
$ws.Cells.Item(12, 40).Value = 17
$ws.Cells.Item(12, 44).Value = 1.95
$ws.Cells.Item(12, 45).Value = 1.85

# Row 15
$ws.Cells.Item(15, 12).Value = 6.5
$ws.Cells.Item(15, 15).Value = 1.18
$ws.Cells.Item(15, 16).Value = 4.5
$ws.Cells.Item(15, 17).Value = 1.58
$ws.Cells.Item(15, 18).Value = 2.25
$ws.Cells.Item(15, 19).Value = 2.5
$ws.Cells.Item(15, 20).Value = 1.5
$ws.Cells.Item(15, 23).Value = 1.83
$ws.Cells.Item(15, 24).Value = 1.83
$ws.Cells.Item(15, 25).Value = 8
$ws.Cells.Item(15, 29).Value = 11
$ws.Cells.Item(15, 30).Value = 23
$ws.Cells.Item(15, 31).Value = 15
$ws.Cells.Item(15, 32).Value = 9
$ws.Cells.Item(15, 33).Value = 19
$ws.Cells.Item(15, 35).Value = 251
$ws.Cells.Item(15, 41).Value = 41

# Row 17
$ws.Cells.Item(17, 7).Value = 4
$ws.Cells.Item(17, 9).Value = 1.75
$ws.Cells.Item(17, 11).Value = 2.6
$ws.Cells.Item(17, 12).Value = 2.25
$ws.Cells.Item(17, 13).Value = 1.01
$ws.Cells.Item(17, 14).Value = 23
$ws.Cells.Item(17, 17).Value = 1.37
$ws.Cells.Item(17, 18).Value = 2.87
$ws.Cells.Item(17, 26).Value = 26
$ws.Cells.Item(17, 27).Value = 13
$ws.Cells.Item(17, 32).Value = 9
$ws.Cells.Item(17, 39).Value = 17

# Row 18
$ws.Cells.Item(18, 7).Value = 5.5
$ws.Cells.Item(18, 9).Value = 1.44
$ws.Cells.Item(18, 11).Value = 2.87
$ws.Cells.Item(18, 12).Value = 1.87
$ws.Cells.Item(18, 14).Value = 29
$ws.Cells.Item(18, 15).Value = 1.08
$ws.Cells.Item(18, 16).Value = 8
$ws.Cells.Item(18, 17).Value = 1.3
$ws.Cells.Item(18, 18).Value = 3.5
$ws.Cells.Item(18, 19).Value = 1.8
$ws.Cells.Item(18, 20).Value = 1.91
$ws.Cells.Item(18, 21).Value = 1.18
$ws.Cells.Item(18, 22).Value = 4.5
$ws.Cells.Item(18, 31).Value = 29
$ws.Cells.Item(18, 33).Value = 15
$ws.Cells.Item(18, 38).Value = 9

# Row 19
$ws.Cells.Item(19, 7).Value = 2.7
$ws.Cells.Item(19, 9).Value = 2.4
$ws.Cells.Item(19, 10).Value = 3.1
$ws.Cells.Item(19, 11).Value = 2.4
$ws.Cells.Item(19, 12).Value = 2.87
$ws.Cells.Item(19, 13).Value = 1.03
$ws.Cells.Item(19, 14).Value = 17
$ws.Cells.Item(19, 15).Value = 1.14
$ws.Cells.Item(19, 16).Value = 5.5
$ws.Cells.Item(19, 17).Value = 1.53
$ws.Cells.Item(19, 18).Value = 2.4
$ws.Cells.Item(19, 19).Value = 2.2
$ws.Cells.Item(19, 20).Value = 1.62
$ws.Cells.Item(19, 21).Value = 1.29
$ws.Cells.Item(19, 22).Value = 3.5
$ws.Cells.Item(19, 23).Value = 1.44
$ws.Cells.Item(19, 24).Value = 2.63
$ws.Cells.Item(19, 31).Value = 17
$ws.Cells.Item(19, 32).Value = 7.5
$ws.Cells.Item(19, 36).Value = 13
$ws.Cells.Item(19, 37).Value = 15
$ws.Cells.Item(19, 39).Value = 23

# Row 20
$ws.Cells.Item(20, 7).Value = 2.63
$ws.Cells.Item(20, 9).Value = 2.75
$ws.Cells.Item(20, 10).Value = 3.2
$ws.Cells.Item(20, 12).Value = 3.25
$ws.Cells.Item(20, 13).Value = 1.05
$ws.Cells.Item(20, 14).Value = 11
$ws.Cells.Item(20, 25).Value = 9.5
$ws.Cells.Item(20, 28).Value = 26
$ws.Cells.Item(20, 33).Value = 12
$ws.Cells.Item(20, 36).Value = 9.5
$ws.Cells.Item(20, 37).Value = 13
$ws.Cells.Item(20, 39).Value = 26
$ws.Cells.Item(20, 42).Value = 1.44
$ws.Cells.Item(20, 43).Value = 2.8

# Row 21
$ws.Cells.Item(21, 17).Value = 2
$ws.Cells.Item(21, 18).Value = 1.85
$ws.Cells.Item(21, 19).Value = 3.25
$ws.Cells.Item(21, 20).Value = 1.33

# Row 22
$ws.Cells.Item(22, 7).Value = 2.05
$ws.Cells.Item(22, 8).Value = 3.25
$ws.Cells.Item(22, 9).Value = 3.75
$ws.Cells.Item(22, 10).Value = 2.75
$ws.Cells.Item(22, 11).Value = 2.05
$ws.Cells.Item(22, 12).Value = 4.33
$ws.Cells.Item(22, 25).Value = 7

# Row 23
$ws.Cells.Item(23, 12).Value = 3
$ws.Cells.Item(23, 23).Value = 2
$ws.Cells.Item(23, 24).Value = 1.73
$ws.Cells.Item(23, 25).Value = 8.5
$ws.Cells.Item(23, 26).Value = 15

# Row 24
$ws.Cells.Item(24, 7).Value = 1.7
$ws.Cells.Item(24, 9).Value = 5
$ws.Cells.Item(24, 14).Value = 8
$ws.Cells.Item(24, 28).Value = 12
$ws.Cells.Item(24, 32).Value = 7.5
$ws.Cells.Item(24, 33).Value = 23
$ws.Cells.Item(24, 42).Value = 1.78
$ws.Cells.Item(24, 43).Value = 2.1

# Row 26
$ws.Cells.Item(26, 7).Value = 3.1
$ws.Cells.Item(26, 9).Value = 2.25
$ws.Cells.Item(26, 10).Value = 3.75
$ws.Cells.Item(26, 11).Value = 2.05
$ws.Cells.Item(26, 12).Value = 3.1
$ws.Cells.Item(26, 13).Value = 1.08
$ws.Cells.Item(26, 14).Value = 8
$ws.Cells.Item(26, 23).Value = 1.83
$ws.Cells.Item(26, 24).Value = 1.83
$ws.Cells.Item(26, 35).Value = 301

# Row 27
$ws.Cells.Item(27, 7).Value = 1.85
$ws.Cells.Item(27, 8).Value = 3.5
$ws.Cells.Item(27, 17).Value = 1.62
$ws.Cells.Item(27, 18).Value = 2.25
$ws.Cells.Item(27, 23).Value = 1.57
$ws.Cells.Item(27, 24).Value = 2.25
$ws.Cells.Item(27, 44).Value = 2.03
$ws.Cells.Item(27, 45).Value = 1.83
